# Update LR-pairs data (Vcan-Tlr2) with recomputed TPM-based statistics.
# Columns E..T for data rows 2..13 are replaced with the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3,1,4.204118999999999,12.612357,0.01983154129720676,0.01983154129720676,2,0.6666666666666666,1.417175333333333,4.251526,0.02388929759887871,0.02388929759887871,5.957973745198,53.621763706782,0.0004737615918934253,0.0004737615918934254),
    @(3,1,4.204118999999999,12.612357,0.01983154129720676,0.01983154129720676,3,1,2.226320666666667,6.678962,0.0375290450698413,0.03752904506984131,9.359717014825998,84.237453133434,0.0007442588071472914,0.0007442588071472917),
    @(3,1,4.204118999999999,12.612357,0.01983154129720676,0.01983154129720676,3,1,55.67910766666667,167.037323,0.93858165733128,0.93858165733128,234.081594444479,2106.734350000311,0.01861352089816604,0.01861352089816604),
    @(3,1,181.0215506666667,543.064652,0.8539093107807857,0.8539093107807858,2,0.6666666666666666,1.417175333333333,4.251526,0.02388929759887871,0.02388929759887871,256.5392764065502,2308.853487658952,0.02039929364769559,0.0203992936476956),
    @(3,1,181.0215506666667,543.064652,0.8539093107807857,0.8539093107807858,3,1,2.226320666666667,6.678962,0.0375290450698413,0.03752904506984131,403.0120193612472,3627.108174251224,0.03204640100984923,0.03204640100984924),
    @(3,1,181.0215506666667,543.064652,0.8539093107807857,0.8539093107807858,3,1,55.67910766666667,167.037323,0.93858165733128,0.93858165733128,10079.11840955629,90712.06568600661,0.8014636161232408,0.8014636161232409),
    @(3,1,26.057747,78.173241,0.1229188055196976,0.1229188055196976,2,0.6666666666666666,1.417175333333333,4.251526,0.02388929759887871,0.02388929759887871,36.92839629064068,332.355566615766,0.002936443925558751,0.002936443925558751),
    @(3,1,26.057747,78.173241,0.1229188055196976,0.1229188055196976,3,1,2.226320666666667,6.678962,0.0375290450698413,0.03752904506984131,58.01290067287134,522.116106055842,0.00461302539227979,0.00461302539227979),
    @(3,1,26.057747,78.173241,0.1229188055196976,0.1229188055196976,3,1,55.67910766666667,167.037323,0.93858165733128,0.93858165733128,1450.872100763761,13057.84890687385,0.1153693362018591,0.1153693362018591),
    @(3,1,0.7081243333333332,2.124373,0.003340342402309973,0.003340342402309974,2,0.6666666666666666,1.417175333333333,4.251526,0.02388929759887871,0.02388929759887871,1.003536338133111,9.031827043198,0.00007979843373093638,0.00007979843373093639),
    @(3,1,0.7081243333333332,2.124373,0.003340342402309973,0.003340342402309974,3,1,2.226320666666667,6.678962,0.0375290450698413,0.03752904506984131,1.576511837869555,14.188606540826,0.000125359860564993,0.000125359860564993),
    @(3,1,0.7081243333333332,2.124373,0.003340342402309973,0.003340342402309974,3,1,55.67910766666667,167.037323,0.93858165733128,0.93858165733128,39.42773099705322,354.849578973479,0.003135184108014044,0.003135184108014044)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $colNum = 5 + $j   # column E = 5
        $ws.Cells.Item($rowNum, $colNum).Value2 = $rowVals[$j]
    }
}
